# Update gh-pages output: add "昆山·星月流光-次元音乐嘉年华" to the 展览 (exhibition)
# sheet (it was originally only on 演出/performance), refresh several
# "想去人数" (interest count) / "最低票价" (min price) figures, and keep the
# combined "全部类型" sheet in sync.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Insert a new row at position 14 for 昆山·星月流光-次元音乐嘉年华 (previously
# listed only under 演出). Existing rows 14-38 shift down to 15-39.
$ws1.Rows.Item(14).Insert()

# Match the bordered/centered index-column formatting used by every other
# row before filling in the new row's data.
$ws1.Cells.Item(13, 1).Copy()
$ws1.Cells.Item(14, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws1.Cells.Item(14, 1).Value = 13
$ws1.Cells.Item(14, 2).NumberFormat = "@"
$ws1.Cells.Item(14, 2).Value = "2024-10-02"
$ws1.Cells.Item(14, 3).Value = "昆山·星月流光-次元音乐嘉年华"
$ws1.Cells.Item(14, 4).Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
$ws1.Cells.Item(14, 5).NumberFormat = "@"
$ws1.Cells.Item(14, 5).Value = "2024.10.02 18:00-10.02 21:00"
$ws1.Cells.Item(14, 6).Value = 1008
$ws1.Cells.Item(14, 7).Value = 19.9
$ws1.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=92044"
$ws1.Cells.Item(14, 9).Value = "//i1.hdslb.com/bfs/openplatform/202409/ZOC3HStE1725591363929.jpeg"

# Column A is a literal running index (0, 1, 2, ...); Insert() shifted the
# rows down but kept each cell's old literal value, so renumber rows 15-39
# back to row-1 to keep the sequence continuous.
for ($r = 15; $r -le 39; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

# Refresh "想去人数" (interest counters) across the sheet.
$ws1.Cells.Item(3, 6).Value = 25
$ws1.Cells.Item(4, 6).Value = 17
$ws1.Cells.Item(5, 6).Value = 15399
$ws1.Cells.Item(9, 6).Value = 15355
$ws1.Cells.Item(11, 6).Value = 8943
$ws1.Cells.Item(12, 6).Value = 364
$ws1.Cells.Item(13, 6).Value = 6
$ws1.Cells.Item(15, 6).Value = 79
$ws1.Cells.Item(16, 6).Value = 193
$ws1.Cells.Item(20, 6).Value = 42
$ws1.Cells.Item(21, 6).Value = 537
$ws1.Cells.Item(24, 6).Value = 59
$ws1.Cells.Item(25, 6).Value = 1105
$ws1.Cells.Item(27, 6).Value = 20
$ws1.Cells.Item(28, 6).Value = 72
$ws1.Cells.Item(34, 6).Value = 240
$ws1.Cells.Item(35, 6).Value = 300
$ws1.Cells.Item(37, 6).Value = 113
$ws1.Cells.Item(38, 6).Value = 5471

# ---------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# 昆山·星月流光-次元音乐嘉年华 now lives solely on 展览, remove it here.
$ws2.Rows.Item(2).Delete()

# Delete() does not renumber the literal index column, fix it up.
$ws2.Cells.Item(2, 1).Value = 1
$ws2.Cells.Item(3, 1).Value = 2

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types combined) - same "想去人数"/"最低票价" refresh,
# no rows move here since it already lists every event.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(3, 6).Value = 25
$ws4.Cells.Item(4, 6).Value = 17
$ws4.Cells.Item(5, 6).Value = 15399
$ws4.Cells.Item(9, 6).Value = 15355
$ws4.Cells.Item(11, 6).Value = 8943
$ws4.Cells.Item(12, 6).Value = 364
$ws4.Cells.Item(13, 6).Value = 6
$ws4.Cells.Item(14, 6).Value = 1008
$ws4.Cells.Item(14, 7).Value = 19.9
$ws4.Cells.Item(15, 6).Value = 79
$ws4.Cells.Item(16, 6).Value = 193
$ws4.Cells.Item(20, 6).Value = 42
$ws4.Cells.Item(21, 6).Value = 537
$ws4.Cells.Item(24, 6).Value = 59
$ws4.Cells.Item(25, 6).Value = 1105
$ws4.Cells.Item(27, 6).Value = 20
$ws4.Cells.Item(28, 6).Value = 72
$ws4.Cells.Item(36, 6).Value = 240
$ws4.Cells.Item(37, 6).Value = 300
$ws4.Cells.Item(39, 6).Value = 113
$ws4.Cells.Item(40, 6).Value = 5471
